$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.339.18'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '1.651.48'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = "'213.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = "'23.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = "'0.0876"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("D12").Value = '1.886.71'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '1.694.51'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("E14").Value = '  -1.65%  '
$ws.Range("D15").Value = "'0.570"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.56%  '
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '27.366.26'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").Value = "'231.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.63%  '
$ws.Range("D19").Value = '0.0₃0724'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("D23").Value = "'9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.03%  '
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("D25").Value = "'147.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("D27").Value = "'15.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.29%  '
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").Value = "'0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("E31").Value = '  -3.92%  '
$ws.Range("E32").Value = '  -1.28%  '
$ws.Range("D33").Value = '1.460.09'
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").Value = "'0.905"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("D38").Value = "'0.571"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  +1.20%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").Value = "'65.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.00%  '
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").Value = '1.794.47'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").Value = "'0.784"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = '  +1.15%  '
$ws.Range("D48").Value = "'88.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("E49").Value = '  -4.04%  '
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  -0.24%  '
